$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.089.24'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.511.08'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  -1.88%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.98'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.89'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -2.90%  '
$ws.Range("B7").Value = 'LidoStakedEther'
$ws.Range("C7").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.498.40'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -2.11%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.613'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -2.87%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  +2.58%  '
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.24'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -2.27%  '
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.069.02'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -2.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.36'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -2.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.513.77'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -1.92%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.009.74'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -1.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.31'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -2.55%  '
$ws.Range("E20").Value = '  -1.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '541.92'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +14.61%  '
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.40'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -4.81%  '
$ws.Range("E24").Value = '  -1.05%  '
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '93.92'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -1.87%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.10'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +1.03%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.93'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -2.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.12'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -1.99%  '
$ws.Range("E30").Value = '  -1.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.26'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -4.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.67'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +4.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '64.54'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -2.77%  '
$ws.Range("E34").Value = '  -3.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '550.41'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -6.64%  '
$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.09'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +8.68%  '
$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '38.07'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -2.12%  '
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +0.24%  '
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.399'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +1.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0765'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -3.96%  '
$ws.Range("E41").Value = '  -2.58%  '
$ws.Range("E42").Value = '  -3.46%  '
$ws.Range("E43").Value = '  -2.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.311.07'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +2.86%  '
$ws.Range("E45").Value = '  -2.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0443'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("E47").Value = '  +2.70%  '
$ws.Range("E48").Value = '  -2.18%  '
$ws.Range("E49").Value = '  -5.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '138.57'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +3.74%  '
